# SSDM-13256: Fixed not correctly behaving compatible with import for vocabularies.
# Changed "Auto generate codes" to "Auto generate code" for samples.
#
# For the EXPERIMENT_TYPE export template, the "Version" column (and its
# example value "1") is removed from the header/example rows, so the
# remaining columns (Code, Description, Validation script) shift left by
# one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (column headers for the type definition): drop "Version", shift
# Code/Description/Validation script one column to the left.
$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "Description"
$ws.Range("C2").Value = "Validation script"
$ws.Range("D2").Clear()

# Row 3 (example values matching the headers above): drop "1", shift the
# remaining example values one column to the left.
$ws.Range("A3").Value = "DEFAULT_EXPERIMENT"
$ws.Range("B3").Value = "Default experiment"
$ws.Range("C3").Value = "test.py"
$ws.Range("D3").Clear()

# Update the stored selection to point at the new first data cell.
$ws.Range("A2").Select()
